# Updated capital structure database
# - Removes the "Financial Services Company SAOG" (row 4) and
#   "Global Financial Investments Holding SAOG" (row 5) records.
# - Refreshes the remaining two Oman "Brokerage & Investment Banking"
#   records (rows 2-3) with new metric values, and drops the now-unused
#   historical_growth_revenue_last_5_years (D) / buybacks_cash_returned (T)
#   figures from row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the two trailing company rows entirely (Financial Services Company SAOG,
# Global Financial Investments Holding SAOG) - shifts dimension to A1:AQ3.
$ws.Rows("4:5").Delete()

# Row 2 (the "3"-ranked company) -------------------------------------------------
$ws.Range("B2").Value = "'1"

# historical_growth_revenue_last_5_years no longer populated for this row.
$ws.Range("D2").ClearContents()

# (quoted "-0" so PowerShell's own numeric-literal parser does not collapse
# the sign before the value reaches the workbook model)
$ws.Range("G2").Value = "-0"
$ws.Range("H2").Value = "-0"
$ws.Range("I2").Value = "-0"
$ws.Range("J2").Value = "-0"
$ws.Range("K2").Value = -6.43
$ws.Range("L2").Value = 32.97435897435897
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = "-0"
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = "-0"

# buybacks_cash_returned no longer populated for this row.
$ws.Range("T2").ClearContents()

$ws.Range("U2").Value = 13.3
$ws.Range("V2").Value = 5.215686274509805
$ws.Range("W2").Value = -0.166580310880829
$ws.Range("X2").Value = 0.04352640958356495
$ws.Range("Y2").Value = -0.2101067204643939
$ws.Range("Z2").Value = -0.05652173913043474
$ws.Range("AB2").Value = 0.04352640958356495
$ws.Range("AC2").Value = -0.04352640958356495
$ws.Range("AD2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -13.3
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 1.237209302325581
$ws.Range("AK2").Value = -0.6650000000000001

# Row 3 (Gulf Investment Services Holding Company SAOG) -------------------------
$ws.Range("K3").Value = -6.43
$ws.Range("L3").Value = 32.97435897435897
$ws.Range("U3").Value = 13.3
$ws.Range("V3").Value = 5.215686274509805
$ws.Range("W3").Value = -0.166580310880829
$ws.Range("X3").Value = 0.04352640958356495
$ws.Range("Y3").Value = -0.2101067204643939
$ws.Range("Z3").Value = -0.05652173913043474
$ws.Range("AB3").Value = 0.04352640958356495
$ws.Range("AC3").Value = -0.04352640958356495
$ws.Range("AD3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -13.3
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 1.237209302325581
$ws.Range("AK3").Value = -0.6650000000000001
